$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure all touched cells keep plain-text formatting (values are text, not numbers)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '278.85'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '6.75%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '27.29'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '0.73%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '4.813'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '2.44%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06285'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '1.01%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.857'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '1.64%'
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = 'MXToken'
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8762'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '2.77%'
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = 'FTXToken'
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9527'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '4.44%'
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1455'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '3.93%'
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.05201'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '9.75%'
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07284'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '2.79%'
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03134'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-0.47%'
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09048'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.13%'
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001560'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '1.09%'
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = 'One'
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0006268'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '1.73%'
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006016'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-1.50%'
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.459'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.54%'
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.272'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '3.04%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.245'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '3.67%'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-0.61%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1310'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-0.08%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.839'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-6.00%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04314'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '1.67%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001172'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-3.48%'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '4.53%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001197'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-0.28%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001685'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '2.74%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04031'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '3.25%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006702'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '62.44%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1153'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '3.69%'
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'LocalTraders'
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/E6DwMU2zXb+localtraders-lct'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.01409'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '1.44%'
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = 'CEJI'
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002098'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-5.12%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005165'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '0.82%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000748'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-0.28%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.329'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '636.28%'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-12.20%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002095'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.28%'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.28%'
